$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1) already exist and are unchanged.

# New data rows (rows 2-7), columns A-T
$data = @(
    @("FAPs","Btc","Erbb3","ECs",1,0.3333333333333333,0.229324,0.687972,0.1345276800989288,0.1345276800989288,2,0.6666666666666666,0.2167736666666666,0.6503209999999999,0.06003183809373913,0.06003183809373914,0.04971140433466666,0.4474026390119999,0.008075943910825223,0.008075943910825226),
    @("FAPs","Btc","Erbb3","FAPs",1,0.3333333333333333,0.229324,0.687972,0.1345276800989288,0.1345276800989288,2,0.6666666666666666,0.3460976666666666,1.038293,0.09584595495126666,0.09584595495126667,0.07936850131066665,0.714316511796,0.01289393396646034,0.01289393396646034),
    @("FAPs","Btc","Erbb3","sCs",1,0.3333333333333333,0.229324,0.687972,0.1345276800989288,0.1345276800989288,3,1,3.048107,9.144321,0.8441222069549942,0.8441222069549942,0.699004089668,6.291036807012,0.1135578022216432,0.1135578022216432),
    @("sCs","Btc","Erbb3","ECs",3,1,1.475336333333334,4.426009000000001,0.8654723199010712,0.8654723199010712,2,0.6666666666666666,0.2167736666666666,0.6503209999999999,0.06003183809373913,0.06003183809373914,0.3198140665432223,2.878326598889,0.05195589418291391,0.05195589418291392),
    @("sCs","Btc","Erbb3","FAPs",3,1,1.475336333333334,4.426009000000001,0.8654723199010712,0.8654723199010712,2,0.6666666666666666,0.3460976666666666,1.038293,0.09584595495126666,0.09584595495126667,0.5106104625152222,4.595494162637,0.08295202098480632,0.08295202098480634),
    @("sCs","Btc","Erbb3","sCs",3,1,1.475336333333334,4.426009000000001,0.8654723199010712,0.8654723199010712,3,1,3.048107,9.144321,0.8441222069549942,0.8441222069549942,4.496983004987667,40.472847044889,0.7305644047333509,0.7305644047333509)
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}
